$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 894, pushing existing rows 894:962 down to 895:963
$ws.Rows.Item(894).Insert()

# Populate the newly inserted row 894 with the new record
$ws.Range("A894").Value = 3
$ws.Range("B894").Value = 'Femacal de La Calera'
$ws.Range("C894").Value = 'Coquimbo'
$ws.Range("D894").Value = 45013
$ws.Range("E894").Value = 5
$ws.Range("F894").Value = 100114001
$ws.Range("G894").Value = 'Papa'
$ws.Range("H894").Value = 'Rosara'
$ws.Range("I894").Value = '1a (cosecha)'
$ws.Range("J894").Value = 510
$ws.Range("K894").Value = 10500
$ws.Range("L894").Value = 11000
$ws.Range("M894").Value = 10775
$ws.Range("N894").Value = '$/saco 25 kilos'
$ws.Range("O894").Value = 'Provincia de Talca'
$ws.Range("P894").Value = 431
$ws.Range("Q894").Value = 25
$ws.Range("R894").Value = 'Hortaliza'
